# case 4 : refresh the simulated curve values in columns A:B (rows 1-32)
# and nudge the column widths to match the recomputed content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 32,2
$arr[0,0] = -0.075632155939885592
$arr[0,1] = 0.075228601081782642
$arr[1,0] = -0.083213360374719159
$arr[1,1] = 0.082574051431584472
$arr[2,0] = -0.093187942962723369
$arr[2,1] = 0.092993391372539591
$arr[3,0] = -0.084993391506566596
$arr[3,1] = 0.084623746673575084
$arr[4,0] = -0.081623746745551173
$arr[4,1] = 0.080372466115737673
$arr[5,0] = 0.018860305832472335
$arr[5,1] = -0.018945040606917374
$arr[6,0] = 0.028945040415603529
$arr[6,1] = -0.028964513430980166
$arr[7,0] = 0.038964513244750698
$arr[7,1] = -0.039038343040094681
$arr[8,0] = 0.041038342977747
$arr[8,1] = -0.041120849876236321
$arr[9,0] = 0.043120849825807994
$arr[9,1] = -0.043125737608711745
$arr[10,0] = 0.00099455094058864546
$arr[10,1] = -0.00099408673406475856
$arr[11,0] = 0.0044940866627483089
$arr[11,1] = -0.0045003199801412208
$arr[12,0] = 0.0080003199109626166
$arr[12,1] = -0.0080112816160404776
$arr[13,0] = 0.016011281484599849
$arr[13,1] = -0.016041805821915389
$arr[14,0] = 0.017041805791581766
$arr[14,1] = -0.017087407573447955
$arr[15,0] = 0.019087407532472955
$arr[15,1] = -0.019290105704546257
$arr[16,0] = 0.021290105672449933
$arr[16,1] = -0.021377727898507892
$arr[17,0] = 0.006669596945126699
$arr[17,1] = -0.0067892754254330612
$arr[18,0] = -0.012090156322797974
$arr[18,1] = 0.012015179461842163
$arr[19,0] = -0.0080151795250653635
$arr[19,1] = 0.0080054032820413568
$arr[20,0] = -0.0040054033460714678
$arr[20,1] = 0.003999999935691001
$arr[21,0] = 0.014613891461907613
$arr[21,1] = -0.014755483848691853
$arr[22,0] = 0.019755483766768833
$arr[22,1] = -0.0200284492189029
$arr[23,0] = -0.020097433986957292
$arr[23,1] = 0.019999999706506344
$arr[24,0] = -0.097232772027654235
$arr[24,1] = 0.097111815061749596
$arr[25,0] = -0.094611815153333723
$arr[25,1] = 0.094454000954023343
$arr[26,0] = -0.091954001052840351
$arr[26,1] = 0.091005761193005519
$arr[27,0] = -0.089005761317936027
$arr[27,1] = 0.088356903051618474
$arr[28,0] = -0.081356903269556469
$arr[28,1] = 0.08116980590382461
$arr[29,0] = -0.021169806880564579
$arr[29,1] = 0.021020452382436883
$arr[30,0] = -0.014020452618883184
$arr[30,1] = 0.014000437864311266
$arr[31,0] = -0.0040004381440255088
$arr[31,1] = 0.0039999998048596552

$ws.Range("A1:B32").Value = $arr

# ColumnWidth snaps to a 1/6-character grid, so these are the nearest
# settable values to the target stored widths (15.7109375 / 16.42578125).
$ws.Columns.Item(1).ColumnWidth = 14.833333333333334
$ws.Columns.Item(2).ColumnWidth = 15.666666666666666
